$d = $word.ActiveDocument

$replacements = @(
    @("79÷7=11, 2", "83÷4=20, 3"),
    @("23÷6=3, 5", "56÷7=8, 0"),
    @("36÷4=9, 0", "33÷5=6, 3"),
    @("56÷8=7, 0", "79÷8=9, 7"),
    @("60÷4=15, 0", "94÷7=13, 3"),
    @("24÷3=8, 0", "32÷3=10, 2"),
    @("96÷9=10, 6", "18÷7=2, 4"),
    @("31÷5=6, 1", "54÷5=10, 4"),
    @("48÷2=24, 0", "18÷2=9, 0"),
    @("14÷5=2, 4", "22÷6=3, 4"),
    @("15÷4=3, 3", "40÷8=5, 0"),
    @("40÷6=6, 4", "70÷5=14, 0"),
    @("29÷4=7, 1", "20÷3=6, 2"),
    @("54÷7=7, 5", "23÷6=3, 5"),
    @("83÷9=9, 2", "29÷6=4, 5"),
    @("26÷4=6, 2", "60÷9=6, 6"),
    @("68÷8=8, 4", "84÷7=12, 0"),
    @("48÷6=8, 0", "11÷2=5, 1"),
    @("41÷3=13, 2", "66÷5=13, 1"),
    @("23÷7=3, 2", "33÷6=5, 3"),
    @("85÷9=9, 4", "90÷7=12, 6"),
    @("18÷8=2, 2", "74÷3=24, 2"),
    @("51÷5=10, 1", "79÷5=15, 4"),
    @("31÷4=7, 3", "40÷2=20, 0"),
    @("62÷9=6, 8", "52÷4=13, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
